$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C slightly to fit the new "Standardabweichung" label
$ws.Columns.Item(3).ColumnWidth = 24.17

# Copy the formatting of the "Durchschnitt" summary row (29) down onto the
# two new summary rows (30 = Median, 31 = Standardabweichung)
$ws.Range("C29:S29").Copy() | Out-Null
$ws.Range("C30:S31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Rows.Item(30).RowHeight = $ws.Rows.Item(29).RowHeight
$ws.Rows.Item(31).RowHeight = $ws.Rows.Item(29).RowHeight

# Row labels
$ws.Range("C30").Value = "Median"
$ws.Range("C31").Value = "Standardabweichung"

# Row 30: Median formulas
$ws.Range("K30").Formula = "=MEDIAN(K7:K25)"
$ws.Range("L30").Formula = "=MEDIAN(L7:L25)"
$ws.Range("N30").Formula = "=MEDIAN(N7:N25)"
$ws.Range("O30").Formula = "=MEDIAN(O7:O25)"
$ws.Range("P30").Formula = "=MEDIAN(P7:P25)"

# Row 31: Standard deviation (population) formulas, rounded to 1 decimal
$ws.Range("K31").Formula = "=ROUND(STDEV.P(K7:K25),1)"
$ws.Range("L31").Formula = "=ROUND(STDEV.P(L7:L25),1)"
$ws.Range("N31").Formula = "=ROUND(STDEV.P(N7:N25),1)"
$ws.Range("O31").Formula = "=ROUND(STDEV.P(O7:O25),1)"
$ws.Range("P31").Formula = "=ROUND(STDEV.P(P7:P25),1)"

# Refresh the view: scroll/selection as recorded by the author at save time
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("N32").Select() | Out-Null
